$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rename variables")

# Update formulas in column P (rows 5-85) to the new format:
#   Table2[[#This Row],[Current]]&"     >>     "&Table2[[#This Row],[New]]
for ($r = 5; $r -le 85; $r++) {
    $cell = $ws.Cells.Item($r, 16)  # column P = 16
    $cell.Formula = "=Table2[[#This Row],[Current]]&`"     >>     `"&Table2[[#This Row],[New]]"
}

# Update the selection shown in the saved sheet view to P5:P85 with active cell P5
$ws.Range("P5:P85").Select()
